$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 998.2308
$ws.Range("J17").Value = 998.2308
$ws.Range("L17").Value = 2994.6924
$ws.Range("N17").Value = -3330.6924
$ws.Range("H43").Value = 4631995.5
$ws.Range("I43").Value = 963.3333
$ws.Range("J43").Value = 6175673
$ws.Range("K43").Value = 963.3333
$ws.Range("L43").Value = 6175673
$ws.Range("M43").Value = -894.3333
$ws.Range("N43").Value = -6175811
$ws.Range("H58").Value = 1679.9474
$ws.Range("J58").Value = 3328
$ws.Range("L58").Value = 9984
$ws.Range("N58").Value = -10284
$ws.Range("H70").Value = 841.9286
$ws.Range("I70").Value = 678.4
$ws.Range("J70").Value = 932.7778
$ws.Range("K70").Value = 2035.2
$ws.Range("L70").Value = 2798.3334
$ws.Range("M70").Value = -1765.2
$ws.Range("N70").Value = -3338.3334
$ws.Range("H73").Value = 841.9286
$ws.Range("I73").Value = 678.4
$ws.Range("J73").Value = 932.7778
$ws.Range("K73").Value = 2035.2
$ws.Range("L73").Value = 2798.3334
$ws.Range("M73").Value = -1099.2
$ws.Range("N73").Value = -4670.3334
$ws.Range("H112").Value = 1768.2778
$ws.Range("I112").Value = 679
$ws.Range("J112").Value = 2312.9167
$ws.Range("K112").Value = 2037
$ws.Range("L112").Value = 6938.750100000001
$ws.Range("M112").Value = -929
$ws.Range("N112").Value = -9154.750100000001
$ws.Range("H132").Value = 6673396
$ws.Range("I132").Value = 7756330
$ws.Range("K132").Value = 23268990
$ws.Range("M132").Value = -23266460
$ws.Range("H135").Value = 31250804
$ws.Range("I135").Value = 546.7143
$ws.Range("J135").Value = 90910380
$ws.Range("K135").Value = 4920.428699999999
$ws.Range("L135").Value = 818193420
$ws.Range("M135").Value = -2385.428699999999
$ws.Range("N135").Value = -818198490
$ws.Range("H141").Value = 798.44446
$ws.Range("I141").Value = 757.17645
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 2271.52935
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = 2908.47065
$ws.Range("N141").Value = -14860
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 934.6
$ws.Range("I2").Value = 999.625
$ws.Range("J2").Value = 674.5
$ws.Range("K2").Value = 999.625
$ws.Range("L2").Value = 674.5
$ws.Range("M2").Value = -886.625
$ws.Range("N2").Value = -900.5
$ws.Range("H15").Value = 3011
$ws.Range("I15").Value = 3011
$ws.Range("K15").Value = 3011
$ws.Range("M15").Value = -2661
$ws.Range("H32").Value = 10255.269
$ws.Range("I32").Value = 7672.319
$ws.Range("J32").Value = 23964.77
$ws.Range("K32").Value = 7672.319
$ws.Range("L32").Value = 23964.77
$ws.Range("M32").Value = -7385.319
$ws.Range("N32").Value = -24538.77
$ws.Range("H45").Value = 1270.1666
$ws.Range("I45").Value = 1005.25
$ws.Range("K45").Value = 1005.25
$ws.Range("M45").Value = -628.25
$ws.Range("H110").Value = 351
$ws.Range("I110").Value = 318.33334
$ws.Range("K110").Value = 318.33334
$ws.Range("M110").Value = 1726.66666
$ws.Range("H116").Value = 934.6
$ws.Range("I116").Value = 999.625
$ws.Range("J116").Value = 674.5
$ws.Range("K116").Value = 999.625
$ws.Range("L116").Value = 674.5
$ws.Range("M116").Value = 1294.375
$ws.Range("N116").Value = -5262.5
$ws.Range("H122").Value = 3157
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3157
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9471
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -14371
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null
$ws.Range("H132").Value = 2743.587
$ws.Range("I132").Value = 2622.5
$ws.Range("K132").Value = 7867.5
$ws.Range("M132").Value = -5337.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 934.6
$ws.Range("I3").Value = 999.625
$ws.Range("J3").Value = 674.5
$ws.Range("K3").Value = 999.625
$ws.Range("L3").Value = 674.5
$ws.Range("M3").Value = -885.625
$ws.Range("N3").Value = -902.5
$ws.Range("H99").Value = 62501196
$ws.Range("I99").Value = 111112080
$ws.Range("J99").Value = 1490
$ws.Range("K99").Value = 111112080
$ws.Range("L99").Value = 1490
$ws.Range("M99").Value = -111110582
$ws.Range("N99").Value = -4486
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1477.4445
$ws.Range("I31").Value = 1382.8223
$ws.Range("J31").Value = 1950.5555
$ws.Range("K31").Value = 1382.8223
$ws.Range("L31").Value = 1950.5555
$ws.Range("M31").Value = -1087.8223
$ws.Range("N31").Value = -2540.5555
$ws.Range("H34").Value = 1477.4445
$ws.Range("I34").Value = 1382.8223
$ws.Range("J34").Value = 1950.5555
$ws.Range("K34").Value = 1382.8223
$ws.Range("L34").Value = 1950.5555
$ws.Range("M34").Value = -1180.8223
$ws.Range("N34").Value = -2354.5555
$ws.Range("H133").Value = 37176
$ws.Range("J133").Value = 37176
$ws.Range("L133").Value = 37176
$ws.Range("N133").Value = -42236
$ws.Range("H141").Value = 690762.9
$ws.Range("J141").Value = 690762.9
$ws.Range("L141").Value = 690762.9
$ws.Range("N141").Value = -701122.9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3511299.8
$ws.Range("J4").Value = 3801499.8
$ws.Range("L4").Value = 11404499.4
$ws.Range("N4").Value = -11404723.4
$ws.Range("H29").Value = 300.84616
$ws.Range("I29").Value = 79.333336
$ws.Range("J29").Value = 367.3
$ws.Range("K29").Value = 238.000008
$ws.Range("L29").Value = 1101.9
$ws.Range("M29").Value = 38.99999199999999
$ws.Range("N29").Value = -1655.9
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = $null
$ws.Range("H92").Value = 510.46155
$ws.Range("I92").Value = 509.55554
$ws.Range("J92").Value = 512.5
$ws.Range("K92").Value = 1528.66662
$ws.Range("L92").Value = 1537.5
$ws.Range("M92").Value = -280.66662
$ws.Range("N92").Value = -4033.5
$ws.Range("H98").Value = 403.83334
$ws.Range("I98").Value = 311.4
$ws.Range("J98").Value = 469.85715
$ws.Range("K98").Value = 934.1999999999999
$ws.Range("L98").Value = 1409.57145
$ws.Range("M98").Value = 563.8000000000001
$ws.Range("N98").Value = -4405.571449999999
$ws.Range("H107").Value = 7218.467
$ws.Range("J107").Value = 8215.154
$ws.Range("L107").Value = 24645.462
$ws.Range("N107").Value = -28485.462
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").Value = $null
$ws.Range("H131").Value = 21309514
$ws.Range("I131").Value = 100000670
$ws.Range("J131").Value = 41632.49
$ws.Range("K131").Value = 300002010
$ws.Range("L131").Value = 124897.47
$ws.Range("M131").Value = -299996970
$ws.Range("N131").Value = -134977.47
$ws.Range("H132").Value = 1332.9166
$ws.Range("I132").Value = 1027.1428
$ws.Range("J132").Value = 1761
$ws.Range("K132").Value = 9244.2852
$ws.Range("L132").Value = 15849
$ws.Range("M132").Value = -6714.2852
$ws.Range("N132").Value = -20909
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4187.375
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4187.375
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 4187.375
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = -6183.375
$ws.Range("H83").Value = 4187.375
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4187.375
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 20936.875
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = -30920.875
$ws.Range("H128").Value = 40000
$ws.Range("I128").Value = 40000
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 40000
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -35020
$ws.Range("N128").Value = $null
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5777.778
$ws.Range("J46").Value = 5777.778
$ws.Range("L46").Value = 5777.778
$ws.Range("N46").Value = -6153.778
$ws.Range("H50").Value = 9056
$ws.Range("J50").Value = 9056
$ws.Range("L50").Value = 9056
$ws.Range("N50").Value = -10330
$ws.Range("H100").Value = 1861.8182
$ws.Range("I100").Value = 1796.6666
$ws.Range("K100").Value = 1796.6666
$ws.Range("M100").Value = -1255.6666
$ws.Range("H132").Value = 3877.8572
$ws.Range("I132").Value = 9999
$ws.Range("J132").Value = 2857.6667
$ws.Range("K132").Value = 29997
$ws.Range("L132").Value = 8573.000100000001
$ws.Range("M132").Value = -27467
$ws.Range("N132").Value = -13633.0001
$ws.Range("H136").Value = 1571.3043
$ws.Range("I136").Value = 1328.3889
$ws.Range("J136").Value = 2445.8
$ws.Range("K136").Value = 3985.1667
$ws.Range("L136").Value = 7337.400000000001
$ws.Range("M136").Value = -1435.1667
$ws.Range("N136").Value = -12437.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -2226
$ws.Range("H24").Value = 154.5
$ws.Range("I24").Value = 154.5
$ws.Range("K24").Value = 154.5
$ws.Range("M24").Value = 75.5
$ws.Range("H44").Value = 10027.333
$ws.Range("J44").Value = 10027.333
$ws.Range("L44").Value = 10027.333
$ws.Range("N44").Value = -11135.333
$ws.Range("H132").Value = 2614.6667
$ws.Range("I132").Value = 2522.2693
$ws.Range("K132").Value = 7566.8079
$ws.Range("M132").Value = -5036.8079
